# Updates Price (D) and Volume(1h) (E) columns of the cryptos list with
# refreshed quotes, matching the commit "Updated cryptos list ...".
# Values that look like plain decimals (e.g. "582.69", "1.00") are forced
# to text via NumberFormat "@" so Excel doesn't silently coerce them to
# numbers (which would also strip the trailing zero / change formatting);
# the style is reset back to "Normal" afterwards so no stray direct
# formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.561.59'
$ws.Range("E2").Value = '  -1.32%  '
$ws.Range("D3").Value = '2.504.99'
$ws.Range("E3").Value = '  -4.90%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.11%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.62'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.95%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").Value = '2.505.16'
$ws.Range("E9").Value = '  -4.85%  '
$ws.Range("E10").Value = '  -0.68%  '
$ws.Range("E11").Value = '  -0.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.33%  '
$ws.Range("E13").Value = '  -2.37%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.74%  '
$ws.Range("D15").Value = '2.970.01'
$ws.Range("E15").Value = '  -4.60%  '
$ws.Range("E16").Value = '  -2.75%  '
$ws.Range("D17").Value = '66.352.60'
$ws.Range("E17").Value = '  -1.61%  '
$ws.Range("D18").Value = '2.505.68'
$ws.Range("E18").Value = '  -4.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.74'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.25'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.24'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.54%  '
$ws.Range("E22").Value = '  -2.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.29%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '69.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = '2.644.28'
$ws.Range("E30").Value = '  -2.98%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '530.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.40%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.98%  '
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("E34").Value = '  -3.00%  '
$ws.Range("E35").Value = '  -4.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("E40").Value = '  +0.38%  '
$ws.Range("E41").Value = '  -3.28%  '
$ws.Range("E42").Value = '  -2.08%  '
$ws.Range("E43").Value = '  -2.25%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.47%  '
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("E46").Value = '  -1.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.35'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.557'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.94%  '
$ws.Range("E49").Value = '  -3.47%  '
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("E51").Value = '  -8.52%  '
